$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 21), matching the date formatting
# already used by the column (copied down from the row above).
$ws.Range("A20").Copy($ws.Range("A21"))
$newDate = Get-Date -Year 2025 -Month 5 -Day 2
$ws.Range("A21").Value = $newDate.Date
$ws.Range("B21").Value = "3 hours"
$ws.Range("C21").Value = "update data and plots"
$ws.Range("E21").Value = "N "

# Update the selection to match the new row, mirroring the author's
# last interaction with the sheet (selecting the newly added row).
$ws.Range("A19:XFD19").Select()
